# Auto-generated edit script: updates market-price / profit cells (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, reflecting a
# refreshed price pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 96.833336
$ws.Range("I2").Value = 96.833336
$ws.Range("K2").Value = 96.833336
$ws.Range("M2").Value = 16.166664

# Row 17
$ws.Range("H17").Value = 6883603.5
$ws.Range("J17").Value = 6883603.5
$ws.Range("L17").Value = 20650810.5
$ws.Range("N17").Value = -20651146.5

# Row 28
$ws.Range("H28").Value = 330.96667
$ws.Range("I28").Value = 187.71428
$ws.Range("K28").Value = 187.71428
$ws.Range("M28").Value = 297.28572

# Row 40
$ws.Range("H40").Value = 5438815.5
$ws.Range("I40").Value = 2669.8823
$ws.Range("J40").Value = 8625521
$ws.Range("K40").Value = 2669.8823
$ws.Range("L40").Value = 8625521
$ws.Range("M40").Value = -2494.8823
$ws.Range("N40").Value = -8625871

# Row 64
$ws.Range("H64").Value = 8187.2856
$ws.Range("I64").Value = 4438.6665
$ws.Range("K64").Value = 4438.6665
$ws.Range("M64").Value = -4190.6665

# Row 67
$ws.Range("H67").Value = 8187.2856
$ws.Range("I67").Value = 4438.6665
$ws.Range("K67").Value = 4438.6665
$ws.Range("M67").Value = -3580.6665

# Row 113
$ws.Range("H113").Value = 27790366
$ws.Range("I113").Value = 31254380
$ws.Range("J113").Value = 78250
$ws.Range("K113").Value = 31254380
$ws.Range("L113").Value = 78250
$ws.Range("M113").Value = -31251126
$ws.Range("N113").Value = -84758

# Row 138
$ws.Range("H138").Value = 5037.269
$ws.Range("J138").Value = 4852.5
$ws.Range("L138").Value = 14557.5
$ws.Range("N138").Value = -24837.5

# Row 141
$ws.Range("H141").Value = 3498.75
$ws.Range("I141").Value = 3498.75
$ws.Range("K141").Value = 10496.25
$ws.Range("M141").Value = -5316.25

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1874.7894
$ws.Range("I2").Value = 1477.8823
$ws.Range("K2").Value = 1477.8823
$ws.Range("M2").Value = -1364.8823

# Row 74
$ws.Range("H74").Value = 3019.3794
$ws.Range("I74").Value = 2010.1666
$ws.Range("K74").Value = 2010.1666
$ws.Range("M74").Value = -1136.1666

# Row 77
$ws.Range("H77").Value = 3019.3794
$ws.Range("I77").Value = 2010.1666
$ws.Range("K77").Value = 10050.833
$ws.Range("M77").Value = -5682.833000000001

# Row 116
$ws.Range("H116").Value = 1874.7894
$ws.Range("I116").Value = 1477.8823
$ws.Range("K116").Value = 1477.8823
$ws.Range("M116").Value = 816.1177

# Row 132
$ws.Range("H132").Value = 2621.1428
$ws.Range("I132").Value = 2624.5
$ws.Range("K132").Value = 7873.5
$ws.Range("M132").Value = -5343.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1874.7894
$ws.Range("I3").Value = 1477.8823
$ws.Range("K3").Value = 1477.8823
$ws.Range("M3").Value = -1363.8823

# Row 107
$ws.Range("H107").Value = 2201.8
$ws.Range("I107").Value = 2331.5293
$ws.Range("K107").Value = 2331.5293
$ws.Range("M107").Value = -411.5293000000001

# Row 134
$ws.Range("H134").Value = 10622.25
$ws.Range("I134").Value = 9995
$ws.Range("K134").Value = 29985
$ws.Range("M134").Value = -27450

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 45460190
$ws.Range("J31").Value = 6276.769
$ws.Range("L31").Value = 6276.769
$ws.Range("N31").Value = -6866.769

# Row 34
$ws.Range("H34").Value = 45460190
$ws.Range("J34").Value = 6276.769
$ws.Range("L34").Value = 6276.769
$ws.Range("N34").Value = -6680.769

# Row 60
$ws.Range("H60").Value = 38048
$ws.Range("I60").Value = 70000
$ws.Range("J60").Value = 30060
$ws.Range("K60").Value = 70000
$ws.Range("L60").Value = 30060
$ws.Range("M60").Value = -69489
$ws.Range("N60").Value = -31082

# Row 99
$ws.Range("H99").Value = 6109.636
$ws.Range("J99").Value = 5599.727
$ws.Range("L99").Value = 5599.727
$ws.Range("N99").Value = -8595.726999999999

# Row 110
$ws.Range("H110").Value = 67439.5
$ws.Range("I110").Value = 50000
$ws.Range("J110").Value = 84879
$ws.Range("K110").Value = 50000
$ws.Range("L110").Value = 84879
$ws.Range("M110").Value = -45910
$ws.Range("N110").Value = -93059

# Row 111
$ws.Range("H111").Value = 80000
$ws.Range("J111").Value = 80000
$ws.Range("L111").Value = 80000
$ws.Range("N111").Value = -88180

# Row 116
$ws.Range("H116").Value = 70000
$ws.Range("J116").Value = 70000
$ws.Range("L116").Value = 70000
$ws.Range("N116").Value = -79178

# Row 126
$ws.Range("H126").Value = 6109.636
$ws.Range("J126").Value = 5599.727
$ws.Range("L126").Value = 16799.181
$ws.Range("N126").Value = -21739.181

# Row 132
$ws.Range("H132").Value = 3953.85
$ws.Range("I132").Value = 2969.9285
$ws.Range("K132").Value = 8909.7855
$ws.Range("M132").Value = -6379.7855

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").Value = ""

# Row 141
$ws.Range("H141").Value = 235294.67
$ws.Range("J141").Value = 267379.47
$ws.Range("L141").Value = 267379.47
$ws.Range("N141").Value = -277739.47

$ws = $wb.Worksheets.Item("CUL")
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 58075
$ws.Range("J32").Value = 58075
$ws.Range("L32").Value = 58075
$ws.Range("N32").Value = -58667

# Row 70
$ws.Range("H70").Value = 8849.951999999999
$ws.Range("I70").Value = 4246
$ws.Range("J70").Value = 11683.154
$ws.Range("K70").Value = 4246
$ws.Range("L70").Value = 11683.154
$ws.Range("M70").Value = -3976
$ws.Range("N70").Value = -12223.154

# Row 73
$ws.Range("H73").Value = 8849.951999999999
$ws.Range("I73").Value = 4246
$ws.Range("J73").Value = 11683.154
$ws.Range("K73").Value = 4246
$ws.Range("L73").Value = 11683.154
$ws.Range("M73").Value = -3310
$ws.Range("N73").Value = -13555.154

# Row 99
$ws.Range("H99").Value = 8724.166999999999
$ws.Range("I99").Value = 5880.909
$ws.Range("J99").Value = 40000
$ws.Range("K99").Value = 5880.909
$ws.Range("L99").Value = 40000
$ws.Range("M99").Value = -3634.909
$ws.Range("N99").Value = -44492

# Row 113
$ws.Range("H113").Value = 573017
$ws.Range("I113").Value = 573017
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 573017
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -570847
$ws.Range("N113").Value = ""

# Row 126
$ws.Range("H126").Value = 5444.2
$ws.Range("I126").Value = 4907
$ws.Range("K126").Value = 14721
$ws.Range("M126").Value = -12251

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2399.5
$ws.Range("I22").Value = 2399.5
$ws.Range("K22").Value = 2399.5
$ws.Range("M22").Value = -2104.5

# Row 27
$ws.Range("H27").Value = 2399.5
$ws.Range("I27").Value = 2399.5
$ws.Range("K27").Value = 2399.5
$ws.Range("M27").Value = -2292.5

# Row 46
$ws.Range("H46").Value = 5686.4443
$ws.Range("J46").Value = 418
$ws.Range("L46").Value = 418
$ws.Range("N46").Value = -794

# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

# Row 59
$ws.Range("H59").Value = 24000
$ws.Range("J59").Value = 24000
$ws.Range("L59").Value = 24000
$ws.Range("N59").Value = -25308

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 10502.667
$ws.Range("J3").Value = 754.5
$ws.Range("L3").Value = 754.5
$ws.Range("N3").Value = -982.5

# Row 62
$ws.Range("H62").Value = 2769.8572
$ws.Range("I62").Value = 1678
$ws.Range("J62").Value = 5499.5
$ws.Range("K62").Value = 1678
$ws.Range("L62").Value = 5499.5
$ws.Range("M62").Value = -1054
$ws.Range("N62").Value = -6747.5

# Row 65
$ws.Range("H65").Value = 2769.8572
$ws.Range("I65").Value = 1678
$ws.Range("J65").Value = 5499.5
$ws.Range("K65").Value = 8390
$ws.Range("L65").Value = 27497.5
$ws.Range("M65").Value = -5270
$ws.Range("N65").Value = -33737.5

# Row 126
$ws.Range("H126").Value = 7050.4517
$ws.Range("I126").Value = 6341.96
$ws.Range("J126").Value = 10002.5
$ws.Range("K126").Value = 19025.88
$ws.Range("L126").Value = 30007.5
$ws.Range("M126").Value = -16555.88
$ws.Range("N126").Value = -34947.5

# Row 132
$ws.Range("H132").Value = 4392.0625
$ws.Range("I132").Value = 4085.96
$ws.Range("J132").Value = 5485.2856
$ws.Range("K132").Value = 12257.88
$ws.Range("L132").Value = 16455.8568
$ws.Range("M132").Value = -9727.880000000001
$ws.Range("N132").Value = -21515.8568

# Row 136
$ws.Range("H136").Value = 4324.233
$ws.Range("I136").Value = 3267.7886
$ws.Range("J136").Value = 6940.1904
$ws.Range("K136").Value = 9803.3658
$ws.Range("L136").Value = 20820.5712
$ws.Range("M136").Value = -7253.3658
$ws.Range("N136").Value = -25920.5712
